$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 6-9 (shrinking the table from 9 rows down to 5)
$ws.Rows("6:9").Delete()

# Update the remaining data rows (A2:B5) with the new LP / Status values
$ws.Range("A2").Value = "LP-049120"
$ws.Range("B2").Value = "Compromisso pendente!"

$ws.Range("A3").Value = "LP-047292"
$ws.Range("B3").Value = "Possui linhas de compra e apontamento!"

$ws.Range("A4").Value = "LP-049245"
$ws.Range("B4").Value = "Compromisso pendente!"

$ws.Range("A5").Value = "LP-049793"
$ws.Range("B5").Value = "Compromisso pendente!"

# Match the new selection shown in the workbook's sheetView
$ws.Range("A2:B5").Select()
